$d = $word.ActiveDocument

# Find the paragraph that contains "kush's edit" so we anchor relative to
# the prior edit rather than assuming a fixed paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*kush's edit*") {
        $target = $p
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

# Move to the very end of that paragraph (just before its paragraph mark)
# and add a blank paragraph followed by a new paragraph with the new text,
# mirroring someone pressing Enter twice and typing "kush's edit 2".
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "kush's edit 2"
